$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.222.17"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "2.347.19"
$ws.Range("E3").Value = "  +3.84%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.649"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0940"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").Value = "2.695.42"
$ws.Range("E13").Value = "  +3.73%  "

$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.837"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").Value = "2.345.53"
$ws.Range("E18").Value = "  +4.42%  "

$ws.Range("D19").Value = "43.157.83"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").Value = "0.0₃0970"
$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("E24").Value = "  +19.98%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.56%  "

$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0685"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "

$ws.Range("E37").Value = "  +7.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("E39").Value = "  -5.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0251"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.56%  "

$ws.Range("E44").Value = "  +7.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0942"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.31%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("D49").Value = "1.433.33"
$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("D50").Value = "2.569.04"
$ws.Range("E50").Value = "  +3.83%  "

$ws.Range("E51").Value = "  -9.28%  "
